{"js": "// Applies the \"Personal info\" bios edit described by the diff:\n//  - Connor paragraph: \"Merriwa\" -> \"Merriwa,\" and \"furbaby\" -> \"furbaby,\"\n//  - Corbin paragraph: replaced with the new first-person bio text\n//  - Two additional blank paragraphs inserted after Corbin's bio (before \"Natalie:\")\n//  - Ollie paragraph: several small wording tweaks\n//  - Trailing bookmark-only paragraph: the stray \"_GoBack\" bookmark is removed\n//    (the paragraph itself stays, now truly empty)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Helper: find a paragraph whose text starts with a given marker.\nfunction findParagraphIndex(items, predicate) {\n  for (let i = 0; i < items.length; i++) {\n    if (predicate(items[i].text)) {\n      return i;\n    }\n  }\n  return -1;\n}\n\nconst items = paragraphs.items;\n\n// --- 1. Connor's paragraph: add the two missing commas ---------------------\nconst connorIdx = findParagraphIndex(items, (t) => t.indexOf(\"S3866963\") === 0);\nconst connorPara = items[connorIdx];\n\nconst merriwaResults = connorPara.search(\"Merriwa\", { matchCase: false });\nmerriwaResults.load(\"items\");\nawait context.sync();\nmerriwaResults.items[0].insertText(\",\", \"After\");\n\nconst furbabyResults = connorPara.search(\"furbaby\", { matchCase: false });\nfurbabyResults.load(\"items\");\nawait context.sync();\nfurbabyResults.items[0].insertText(\",\", \"After\");\n\nawait context.sync();\n\n// --- 2. Corbin's paragraph: replace with the new bio ------------------------\nconst corbinIdx = findParagraphIndex(\n  items,\n  (t) => t.indexOf(\"s3855159\") === 0\n);\nconst corbinPara = items[corbinIdx];\nconst newCorbinText =\n  \"My name is Corbin, RMIT ID: s3855159, from team XVI. Hailing from Melbourne City and originally from Country Victoria, I now live on the sunny Mornington Peninsula. My day job is in hospitality and tourism management but I\\u2019m always pursuing a new side-hustle. I\\u2019m an enormous music lover and spend most of my free time singing and playing guitar, reading a good book or involving myself in social or philosophical discussions. I\\u2019ve been a gamer for my whole life, and I don\\u2019t really remember a time when I haven\\u2019t owned some sort of gaming console. My interest in IT was spurred when I took a short course in Python. Tech had always interested me but learning a small amount of code showed me that it was something I could learn and not as out-of-reach as it appeared. I\\u2019m particularly interested in artificial intelligence and the future of computing.\";\ncorbinPara.insertText(newCorbinText, \"Replace\");\nawait context.sync();\n\n// --- 3. Insert two extra blank paragraphs right after Corbin's bio ---------\n// (the blank line that already follows Corbin's bio becomes the first of three)\ncorbinPara.insertParagraph(\"\", \"After\");\ncorbinPara.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// --- 4. Ollie's paragraph: apply the wording tweaks -------------------------\nparagraphs.load(\"items/text\");\nawait context.sync();\nconst items2 = paragraphs.items;\nconst ollieIdx = findParagraphIndex(\n  items2,\n  (t) => t.indexOf(\"S3861675\") === 0\n);\nconst olliePara = items2[ollieIdx];\nconst newOllieText =\n  \"S3861675. My name is Oliver, I am 16 and was born in Australia. I enjoy playing video games, watching shows, going out with friends and using software such as Unity to explore cool ideas. I have always enjoyed using technology, whether it was creating my own retro arcades with raspberry pi\\u2019s or small fun games in Unity to mess around in, with friends. I have never made a game with a serious intent to either sell or release it, but I have made many to share with friends and play together for the next week seeing who can get the highest score. IT leaves almost no limit to creativity and that\\u2019s why I like it so much. I would love to get a job as a game developer in a company, but it has also been my goal to work either by myself or in a small team just having fun whether it is in YouTube or making Indie Titles.\";\nolliePara.insertText(newOllieText, \"Replace\");\nawait context.sync();\n\n// --- 5. Remove the stray \"_GoBack\" bookmark on the trailing paragraph ------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Applies the \"Personal info\" bios edit described by the diff:\n#  - Connor paragraph: \"Merriwa\" -> \"Merriwa,\" and \"furbaby\" -> \"furbaby,\"\n#  - Corbin paragraph: replaced with the new first-person bio text\n#  - Two additional blank paragraphs inserted after Corbin's bio (before \"Natalie:\")\n#  - Ollie paragraph: several small wording tweaks\n#  - Trailing bookmark-only paragraph: the stray \"_GoBack\" bookmark is removed\n#    (the paragraph itself stays, now truly empty)\n\n$d = $word.ActiveDocument\n\n# --- 1. Connor's paragraph: add the two missing commas ---------------------\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*S3866963*\") {\n    $r = $p.Range\n    $r.End = $r.End - 1\n    $newText = $r.Text\n    $newText = $newText -replace \"Merriwa NSW\", \"Merriwa, NSW\"\n    $newText = $newText -replace \"furbaby Turbo\", \"furbaby, Turbo\"\n    $r.Text = $newText\n    break\n  }\n}\n\n# --- 2. Corbin's paragraph: replace with the new bio ------------------------\n$corbinIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*s3855159*\") {\n    $corbinIndex = $i\n    $r = $p.Range\n    $r.End = $r.End - 1\n    $r.Text = \"My name is Corbin, RMIT ID: s3855159, from team XVI. Hailing from Melbourne City and originally from Country Victoria, I now live on the sunny Mornington Peninsula. My day job is in hospitality and tourism management but I\u2019m always pursuing a new side-hustle. I\u2019m an enormous music lover and spend most of my free time singing and playing guitar, reading a good book or involving myself in social or philosophical discussions. I\u2019ve been a gamer for my whole life, and I don\u2019t really remember a time when I haven\u2019t owned some sort of gaming console. My interest in IT was spurred when I took a short course in Python. Tech had always interested me but learning a small amount of code showed me that it was something I could learn and not as out-of-reach as it appeared. I\u2019m particularly interested in artificial intelligence and the future of computing.\"\n    break\n  }\n}\n\n# --- 3. Insert two extra blank paragraphs right after Corbin's bio ---------\n# (the blank line that already follows Corbin's bio becomes the first of three)\n$corbinPara = $d.Paragraphs.Item($corbinIndex)\n$endOfCorbin = $corbinPara.Range.End\n$insertRange = $d.Range($endOfCorbin, $endOfCorbin)\n$insertRange.InsertParagraphAfter()\n$insertRange.InsertParagraphAfter()\n\n# --- 4. Ollie's paragraph: apply the wording tweaks -------------------------\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*S3861675*\") {\n    $r = $p.Range\n    $r.End = $r.End - 1\n    $r.Text = \"S3861675. My name is Oliver, I am 16 and was born in Australia. I enjoy playing video games, watching shows, going out with friends and using software such as Unity to explore cool ideas. I have always enjoyed using technology, whether it was creating my own retro arcades with raspberry pi\u2019s or small fun games in Unity to mess around in, with friends. I have never made a game with a serious intent to either sell or release it, but I have made many to share with friends and play together for the next week seeing who can get the highest score. IT leaves almost no limit to creativity and that\u2019s why I like it so much. I would love to get a job as a game developer in a company, but it has also been my goal to work either by myself or in a small team just having fun whether it is in YouTube or making Indie Titles.\"\n    break\n  }\n}\n\n# --- 5. Remove the stray \"_GoBack\" bookmark on the trailing paragraph ------\n$d.Bookmarks.ShowHidden = $true\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n"}
